# Saldo_guide.xlsx update
#  - sheet renamed from the 2024-11-13 run to the 2024-11-14 run
#  - every "Data" (column G) value bumped one day forward (45609 -> 45610)
#  - active cell / selection left on K13

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet to match the new export timestamp.
$ws.Name = "IClientBalance-20241114-082649-"

# Bump every date in column G (rows 2-274) from 45609 to 45610 (13th -> 14th Nov 2024).
$ws.Range("G2:G274").Value = 45610

# Move the active selection to K13, as in the edited workbook.
$ws.Range("K13").Select()
